# Apply cryptos.xlsx data refresh updates
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

    $ws.Range('D2').NumberFormat = '@'
    $ws.Range('D2').Value = '54.253.39'
    $ws.Range('D2').Style = 'Normal'
    $ws.Range('E2').NumberFormat = '@'
    $ws.Range('E2').Value = '  -6.09%  '
    $ws.Range('E2').Style = 'Normal'
    $ws.Range('D3').NumberFormat = '@'
    $ws.Range('D3').Value = '2.860.92'
    $ws.Range('D3').Style = 'Normal'
    $ws.Range('E3').NumberFormat = '@'
    $ws.Range('E3').Value = '  -9.74%  '
    $ws.Range('E3').Style = 'Normal'
    $ws.Range('E4').NumberFormat = '@'
    $ws.Range('E4').Value = '  +0.10%  '
    $ws.Range('E4').Style = 'Normal'
    $ws.Range('D5').NumberFormat = '@'
    $ws.Range('D5').Value = '465.64'
    $ws.Range('D5').Style = 'Normal'
    $ws.Range('E5').NumberFormat = '@'
    $ws.Range('E5').Value = '  -12.35%  '
    $ws.Range('E5').Style = 'Normal'
    $ws.Range('D6').NumberFormat = '@'
    $ws.Range('D6').Value = '123.40'
    $ws.Range('D6').Style = 'Normal'
    $ws.Range('E6').NumberFormat = '@'
    $ws.Range('E6').Value = '  -8.37%  '
    $ws.Range('E6').Style = 'Normal'
    $ws.Range('E7').NumberFormat = '@'
    $ws.Range('E7').Value = '  -0.05%  '
    $ws.Range('E7').Style = 'Normal'
    $ws.Range('D8').NumberFormat = '@'
    $ws.Range('D8').Value = '2.861.93'
    $ws.Range('D8').Style = 'Normal'
    $ws.Range('E8').NumberFormat = '@'
    $ws.Range('E8').Value = '  -9.71%  '
    $ws.Range('E8').Style = 'Normal'
    $ws.Range('D9').NumberFormat = '@'
    $ws.Range('D9').Value = '0.399'
    $ws.Range('D9').Style = 'Normal'
    $ws.Range('E9').NumberFormat = '@'
    $ws.Range('E9').Value = '  -11.78%  '
    $ws.Range('E9').Style = 'Normal'
    $ws.Range('D10').NumberFormat = '@'
    $ws.Range('D10').Value = '6.55'
    $ws.Range('D10').Style = 'Normal'
    $ws.Range('E10').NumberFormat = '@'
    $ws.Range('E10').Value = '  -9.90%  '
    $ws.Range('E10').Style = 'Normal'
    $ws.Range('D11').NumberFormat = '@'
    $ws.Range('D11').Value = '0.0944'
    $ws.Range('D11').Style = 'Normal'
    $ws.Range('E11').NumberFormat = '@'
    $ws.Range('E11').Value = '  -15.16%  '
    $ws.Range('E11').Style = 'Normal'
    $ws.Range('E12').NumberFormat = '@'
    $ws.Range('E12').Value = '  -18.40%  '
    $ws.Range('E12').Style = 'Normal'
    $ws.Range('E13').NumberFormat = '@'
    $ws.Range('E13').Value = '  -5.06%  '
    $ws.Range('E13').Style = 'Normal'
    $ws.Range('D14').NumberFormat = '@'
    $ws.Range('D14').Value = '3.353.91'
    $ws.Range('D14').Style = 'Normal'
    $ws.Range('E14').NumberFormat = '@'
    $ws.Range('E14').Value = '  -9.73%  '
    $ws.Range('E14').Style = 'Normal'
    $ws.Range('D15').NumberFormat = '@'
    $ws.Range('D15').Value = '22.86'
    $ws.Range('D15').Style = 'Normal'
    $ws.Range('E15').NumberFormat = '@'
    $ws.Range('E15').Value = '  -11.37%  '
    $ws.Range('E15').Style = 'Normal'
    $ws.Range('D16').NumberFormat = '@'
    $ws.Range('D16').Value = '54.238.15'
    $ws.Range('D16').Style = 'Normal'
    $ws.Range('E16').NumberFormat = '@'
    $ws.Range('E16').Value = '  -6.27%  '
    $ws.Range('E16').Style = 'Normal'
    $ws.Range('D17').NumberFormat = '@'
    $ws.Range('D17').Value = '2.865.71'
    $ws.Range('D17').Style = 'Normal'
    $ws.Range('E17').NumberFormat = '@'
    $ws.Range('E17').Value = '  -9.78%  '
    $ws.Range('E17').Style = 'Normal'
    $ws.Range('D18').NumberFormat = '@'
    $ws.Range('D18').Value = '0.0000130'
    $ws.Range('D18').Style = 'Normal'
    $ws.Range('E18').NumberFormat = '@'
    $ws.Range('E18').Value = '  -15.66%  '
    $ws.Range('E18').Style = 'Normal'
    $ws.Range('D19').NumberFormat = '@'
    $ws.Range('D19').Value = '5.24'
    $ws.Range('D19').Style = 'Normal'
    $ws.Range('E19').NumberFormat = '@'
    $ws.Range('E19').Value = '  -10.14%  '
    $ws.Range('E19').Style = 'Normal'
    $ws.Range('D20').NumberFormat = '@'
    $ws.Range('D20').Value = '11.18'
    $ws.Range('D20').Style = 'Normal'
    $ws.Range('E20').NumberFormat = '@'
    $ws.Range('E20').Value = '  -15.60%  '
    $ws.Range('E20').Style = 'Normal'
    $ws.Range('D21').NumberFormat = '@'
    $ws.Range('D21').Value = '6.94'
    $ws.Range('D21').Style = 'Normal'
    $ws.Range('E21').NumberFormat = '@'
    $ws.Range('E21').Value = '  -14.14%  '
    $ws.Range('E21').Style = 'Normal'
    $ws.Range('D22').NumberFormat = '@'
    $ws.Range('D22').Value = '291.10'
    $ws.Range('D22').Style = 'Normal'
    $ws.Range('E22').NumberFormat = '@'
    $ws.Range('E22').Value = '  -18.64%  '
    $ws.Range('E22').Style = 'Normal'
    $ws.Range('D23').NumberFormat = '@'
    $ws.Range('D23').Value = '0.998'
    $ws.Range('D23').Style = 'Normal'
    $ws.Range('E23').NumberFormat = '@'
    $ws.Range('E23').Value = '  -0.06%  '
    $ws.Range('E23').Style = 'Normal'
    $ws.Range('E24').NumberFormat = '@'
    $ws.Range('E24').Value = '  -16.33%  '
    $ws.Range('E24').Style = 'Normal'
    $ws.Range('D25').NumberFormat = '@'
    $ws.Range('D25').Value = '57.75'
    $ws.Range('D25').Style = 'Normal'
    $ws.Range('E25').NumberFormat = '@'
    $ws.Range('E25').Value = '  -16.98%  '
    $ws.Range('E25').Style = 'Normal'
    $ws.Range('E26').NumberFormat = '@'
    $ws.Range('E26').Value = '  +0.33%  '
    $ws.Range('E26').Style = 'Normal'
    $ws.Range('E27').NumberFormat = '@'
    $ws.Range('E27').Value = '  +0.09%  '
    $ws.Range('E27').Style = 'Normal'
    $ws.Range('D28').NumberFormat = '@'
    $ws.Range('D28').Value = '0.148'
    $ws.Range('D28').Style = 'Normal'
    $ws.Range('E28').NumberFormat = '@'
    $ws.Range('E28').Value = '  -12.02%  '
    $ws.Range('E28').Style = 'Normal'
    $ws.Range('D29').NumberFormat = '@'
    $ws.Range('D29').Value = '0.0₃0782'
    $ws.Range('D29').Style = 'Normal'
    $ws.Range('E29').NumberFormat = '@'
    $ws.Range('E29').Value = '  -17.90%  '
    $ws.Range('E29').Style = 'Normal'
    $ws.Range('E30').NumberFormat = '@'
    $ws.Range('E30').Value = '  -13.68%  '
    $ws.Range('E30').Style = 'Normal'
    $ws.Range('D31').NumberFormat = '@'
    $ws.Range('D31').Value = '6.02'
    $ws.Range('D31').Style = 'Normal'
    $ws.Range('E31').NumberFormat = '@'
    $ws.Range('E31').Value = '  -13.72%  '
    $ws.Range('E31').Style = 'Normal'
    $ws.Range('E32').NumberFormat = '@'
    $ws.Range('E32').Value = '  -10.50%  '
    $ws.Range('E32').Style = 'Normal'
    $ws.Range('E33').NumberFormat = '@'
    $ws.Range('E33').Value = '  -16.47%  '
    $ws.Range('E33').Style = 'Normal'
    $ws.Range('D34').NumberFormat = '@'
    $ws.Range('D34').Value = '18.24'
    $ws.Range('D34').Style = 'Normal'
    $ws.Range('E34').NumberFormat = '@'
    $ws.Range('E34').Value = '  -15.77%  '
    $ws.Range('E34').Style = 'Normal'
    $ws.Range('D35').NumberFormat = '@'
    $ws.Range('D35').Value = '136.84'
    $ws.Range('D35').Style = 'Normal'
    $ws.Range('E35').NumberFormat = '@'
    $ws.Range('E35').Value = '  -14.55%  '
    $ws.Range('E35').Style = 'Normal'
    $ws.Range('E36').NumberFormat = '@'
    $ws.Range('E36').Value = '  -18.03%  '
    $ws.Range('E36').Style = 'Normal'
    $ws.Range('E37').NumberFormat = '@'
    $ws.Range('E37').Value = '  -15.61%  '
    $ws.Range('E37').Style = 'Normal'
    $ws.Range('E38').NumberFormat = '@'
    $ws.Range('E38').Value = '  -16.92%  '
    $ws.Range('E38').Style = 'Normal'
    $ws.Range('D39').NumberFormat = '@'
    $ws.Range('D39').Value = '22.63'
    $ws.Range('D39').Style = 'Normal'
    $ws.Range('E39').NumberFormat = '@'
    $ws.Range('E39').Value = '  -12.82%  '
    $ws.Range('E39').Style = 'Normal'
    $ws.Range('D40').NumberFormat = '@'
    $ws.Range('D40').Value = '2.888.43'
    $ws.Range('D40').Style = 'Normal'
    $ws.Range('E40').NumberFormat = '@'
    $ws.Range('E40').Value = '  -9.86%  '
    $ws.Range('E40').Style = 'Normal'
    $ws.Range('E41').NumberFormat = '@'
    $ws.Range('E41').Value = '  -0.01%  '
    $ws.Range('E41').Style = 'Normal'
    $ws.Range('D42').NumberFormat = '@'
    $ws.Range('D42').Value = '0.0606'
    $ws.Range('D42').Style = 'Normal'
    $ws.Range('E42').NumberFormat = '@'
    $ws.Range('E42').Value = '  -13.96%  '
    $ws.Range('E42').Style = 'Normal'
    $ws.Range('E43').NumberFormat = '@'
    $ws.Range('E43').Value = '  -13.03%  '
    $ws.Range('E43').Style = 'Normal'
    $ws.Range('E44').NumberFormat = '@'
    $ws.Range('E44').Value = '  -15.08%  '
    $ws.Range('E44').Style = 'Normal'
    $ws.Range('D45').NumberFormat = '@'
    $ws.Range('D45').Value = '0.918'
    $ws.Range('D45').Style = 'Normal'
    $ws.Range('E45').NumberFormat = '@'
    $ws.Range('E45').Value = '  -15.88%  '
    $ws.Range('E45').Style = 'Normal'
    $ws.Range('D46').NumberFormat = '@'
    $ws.Range('D46').Value = '1.28'
    $ws.Range('D46').Style = 'Normal'
    $ws.Range('E46').NumberFormat = '@'
    $ws.Range('E46').Value = '  -13.39%  '
    $ws.Range('E46').Style = 'Normal'
    $ws.Range('B47').Value = 'Filecoin'
    $ws.Range('C47').Value = 'https://coinranking.com/coin/ymQub4fuB+filecoin-fil'
    $ws.Range('D47').NumberFormat = '@'
    $ws.Range('D47').Value = '3.31'
    $ws.Range('D47').Style = 'Normal'
    $ws.Range('E47').NumberFormat = '@'
    $ws.Range('E47').Value = '  -17.17%  '
    $ws.Range('E47').Style = 'Normal'
    $ws.Range('B48').Value = 'Maker'
    $ws.Range('C48').Value = 'https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr'
    $ws.Range('D48').NumberFormat = '@'
    $ws.Range('D48').Value = '2.020.31'
    $ws.Range('D48').Style = 'Normal'
    $ws.Range('E48').NumberFormat = '@'
    $ws.Range('E48').Value = '  -11.17%  '
    $ws.Range('E48').Style = 'Normal'
    $ws.Range('D49').NumberFormat = '@'
    $ws.Range('D49').Value = '5.26'
    $ws.Range('D49').Style = 'Normal'
    $ws.Range('E49').NumberFormat = '@'
    $ws.Range('E49').Value = '  -15.48%  '
    $ws.Range('E49').Style = 'Normal'
    $ws.Range('E50').NumberFormat = '@'
    $ws.Range('E50').Value = '  -11.36%  '
    $ws.Range('E50').Style = 'Normal'
    $ws.Range('D51').NumberFormat = '@'
    $ws.Range('D51').Value = '17.41'
    $ws.Range('D51').Style = 'Normal'
    $ws.Range('E51').NumberFormat = '@'
    $ws.Range('E51').Value = '  -15.91%  '
    $ws.Range('E51').Style = 'Normal'
